$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet and add the new history sheet ---
$ws1 = $wb.ActiveSheet
$ws1.Name = "Stato Attuale"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Storico Passaggi"

# --- 2. Update the current-state sheet: GX666SK moved from PASQUARELLI to FINE RENT ---
$targaRange = $ws1.Range("A101")
if ($targaRange.Value -ne "GX666SK") {
    # Fallback: locate the row holding GX666SK if the layout ever shifts
    $found = $ws1.Cells.Find("GX666SK")
    if ($found -ne $null) {
        $targaRow = $found.Row
    } else {
        $targaRow = 101
    }
} else {
    $targaRow = 101
}

$ws1.Cells.Item($targaRow, 2).Value = "FINE RENT"

$dateCell = $ws1.Cells.Item($targaRow, 3)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-18"
$dateCell.ClearFormats()

# --- 3. Build the history sheet header ---
$ws2.Range("A1").Value = "Targa"
$ws2.Range("B1").Value = "Operatore_Precedente"
$ws2.Range("C1").Value = "Nuovo_Operatore"
$ws2.Range("D1").Value = "Data_Cambio"

$headerRange = $ws2.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- 4. Log the change as a history row ---
$ws2.Range("A2").Value = "GX666SK"
$ws2.Range("B2").Value = "PASQUARELLI"
$ws2.Range("C2").Value = "FINE RENT"

$histDateCell = $ws2.Range("D2")
$histDateCell.NumberFormat = "@"
$histDateCell.Value = "2025-12-18"
$histDateCell.ClearFormats()
